$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 32576.715
$ws.Range("I19").Value = 916
$ws.Range("J19").Value = 111728.5
$ws.Range("K19").Value = 916
$ws.Range("L19").Value = 111728.5
$ws.Range("M19").Value = -741
$ws.Range("N19").Value = -112078.5
$ws.Range("H33").Value = 6946847
$ws.Range("I33").Value = 10000234
$ws.Range("J33").Value = 7330.909
$ws.Range("K33").Value = 10000234
$ws.Range("L33").Value = 7330.909
$ws.Range("M33").Value = -10000005
$ws.Range("N33").Value = -7788.909
$ws.Range("H40").Value = 875
$ws.Range("J40").Value = 1150
$ws.Range("L40").Value = 1150
$ws.Range("N40").Value = -1500
$ws.Range("H43").Value = 2365.923
$ws.Range("I43").Value = 3680
$ws.Range("J43").Value = 1544.625
$ws.Range("K43").Value = 3680
$ws.Range("L43").Value = 1544.625
$ws.Range("M43").Value = -3611
$ws.Range("N43").Value = -1682.625
$ws.Range("H98").Value = 4698.8965
$ws.Range("I98").Value = 4688.143
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 4688.143
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = -3190.143
$ws.Range("N98").Value = -7996
$ws.Range("H106").Value = 4472.3184
$ws.Range("I106").Value = 2352.2632
$ws.Range("J106").Value = 17899.334
$ws.Range("K106").Value = 2352.2632
$ws.Range("L106").Value = 17899.334
$ws.Range("M106").Value = -1721.2632
$ws.Range("N106").Value = -19161.334
$ws.Range("H116").Value = 3296462.2
$ws.Range("I116").Value = 4237051.5
$ws.Range("J116").Value = 4399.1665
$ws.Range("K116").Value = 4237051.5
$ws.Range("L116").Value = 4399.1665
$ws.Range("M116").Value = -4233609.5
$ws.Range("N116").Value = -11283.1665
$ws.Range("H122").Value = 4698.8965
$ws.Range("I122").Value = 4688.143
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14064.429
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -11614.429
$ws.Range("N122").Value = -19900
$ws.Range("H125").Value = 5404.8
$ws.Range("I125").Value = 5251.4443
$ws.Range("J125").Value = 5530.273
$ws.Range("K125").Value = 47262.9987
$ws.Range("L125").Value = 49772.457
$ws.Range("M125").Value = -44802.9987
$ws.Range("N125").Value = -54692.457

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3545.5398
$ws.Range("I32").Value = 2404.7068
$ws.Range("J32").Value = 16779.2
$ws.Range("K32").Value = 2404.7068
$ws.Range("L32").Value = 16779.2
$ws.Range("M32").Value = -2117.7068
$ws.Range("N32").Value = -17353.2
$ws.Range("H74").Value = 15769.448
$ws.Range("I74").Value = 1694.8636
$ws.Range("J74").Value = 60003.855
$ws.Range("K74").Value = 1694.8636
$ws.Range("L74").Value = 60003.855
$ws.Range("M74").Value = -820.8635999999999
$ws.Range("N74").Value = -61751.855
$ws.Range("H77").Value = 15769.448
$ws.Range("I77").Value = 1694.8636
$ws.Range("J77").Value = 60003.855
$ws.Range("K77").Value = 8474.317999999999
$ws.Range("L77").Value = 300019.275
$ws.Range("M77").Value = -4106.317999999999
$ws.Range("N77").Value = -308755.275
$ws.Range("H122").Value = 2023.2084
$ws.Range("I122").Value = 1546.1875
$ws.Range("J122").Value = 2977.25
$ws.Range("K122").Value = 4638.5625
$ws.Range("L122").Value = 8931.75
$ws.Range("M122").Value = -2188.5625
$ws.Range("N122").Value = -13831.75
$ws.Range("H132").Value = 4015937.2
$ws.Range("I132").Value = 2225.2856
$ws.Range("J132").Value = 11039933
$ws.Range("K132").Value = 6675.8568
$ws.Range("L132").Value = 33119799
$ws.Range("M132").Value = -4145.8568
$ws.Range("N132").Value = -33124859

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9976.440000000001
$ws.Range("I99").Value = 12444.895
$ws.Range("K99").Value = 12444.895
$ws.Range("M99").Value = -10946.895
$ws.Range("H134").Value = 17916.062
$ws.Range("I134").Value = 14047.218
$ws.Range("J134").Value = 27803.111
$ws.Range("K134").Value = 42141.654
$ws.Range("L134").Value = 83409.333
$ws.Range("M134").Value = -39606.654
$ws.Range("N134").Value = -88479.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11652.779
$ws.Range("I31").Value = 3903.423
$ws.Range("J31").Value = 16450
$ws.Range("K31").Value = 3903.423
$ws.Range("L31").Value = 16450
$ws.Range("M31").Value = -3608.423
$ws.Range("N31").Value = -17040
$ws.Range("H34").Value = 11652.779
$ws.Range("I34").Value = 3903.423
$ws.Range("J34").Value = 16450
$ws.Range("K34").Value = 3903.423
$ws.Range("L34").Value = 16450
$ws.Range("M34").Value = -3701.423
$ws.Range("N34").Value = -16854
$ws.Range("H86").Value = 10704.1875
$ws.Range("I86").Value = 11847.7
$ws.Range("J86").Value = 8798.333000000001
$ws.Range("K86").Value = 11847.7
$ws.Range("L86").Value = 8798.333000000001
$ws.Range("M86").Value = -10724.7
$ws.Range("N86").Value = -11044.333
$ws.Range("H89").Value = 10704.1875
$ws.Range("I89").Value = 11847.7
$ws.Range("J89").Value = 8798.333000000001
$ws.Range("K89").Value = 59238.5
$ws.Range("L89").Value = 43991.665
$ws.Range("M89").Value = -53622.5
$ws.Range("N89").Value = -55223.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5124.6665
$ws.Range("J68").Value = 6249.5
$ws.Range("L68").Value = 18748.5
$ws.Range("N68").Value = -20370.5
$ws.Range("H71").Value = 5124.6665
$ws.Range("J71").Value = 6249.5
$ws.Range("L71").Value = 56245.5
$ws.Range("N71").Value = -64357.5
$ws.Range("H92").Value = 218.88235
$ws.Range("I92").Value = 165.125
$ws.Range("J92").Value = 266.66666
$ws.Range("K92").Value = 495.375
$ws.Range("L92").Value = 799.9999799999999
$ws.Range("M92").Value = 752.625
$ws.Range("N92").Value = -3295.99998
$ws.Range("H122").Value = 14351333
$ws.Range("J122").Value = 3548623
$ws.Range("L122").Value = 31937607
$ws.Range("N122").Value = -31942507
$ws.Range("H137").Value = 5504.769
$ws.Range("I137").Value = 1755.8
$ws.Range("J137").Value = 6397.381
$ws.Range("K137").Value = 5267.4
$ws.Range("L137").Value = 19192.143
$ws.Range("M137").Value = -167.3999999999996
$ws.Range("N137").Value = -29392.143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 43329.668
$ws.Range("I122").Value = 39995
$ws.Range("J122").Value = 49999
$ws.Range("K122").Value = 119985
$ws.Range("L122").Value = 149997
$ws.Range("M122").Value = -117535
$ws.Range("N122").Value = -154897

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 8638.23
$ws.Range("I100").Value = 5787.25
$ws.Range("K100").Value = 5787.25
$ws.Range("M100").Value = -5246.25
